$d = $word.ActiveDocument

# --- Locate the date paragraph text ("12 October 2014") ----------------
$dateRange = $d.Content.Duplicate
$found = $dateRange.Find.Execute("12 October 2014", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the date text '12 October 2014' in the document."
}

$dateStart = $dateRange.Start

# --- Remove the existing "_GoBack" bookmark wherever it currently lives -
# The author's edit relocates this bookmark from the end of the abstract
# paragraph to the middle of the date ("Octo|ber"). Bookmark names must be
# unique in a document, so drop the old one before adding the new one.
try {
    $oldBm = $d.Bookmarks("_GoBack")
    $oldBm.Delete()
} catch {
    # no existing "_GoBack" bookmark - nothing to remove
}

# --- "12 October 2014" -> "11 October 2014" -----------------------------
$r12 = $d.Range($dateStart, $dateStart + 2)
$r12.Text = "11"

# --- Re-insert the "_GoBack" bookmark between "Octo" and "ber" ---------
# "11 Octo" is 7 characters from the start of the (now updated) date text.
$bmPos = $dateStart + 7
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "done"
